$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the team member name in A2
$ws.Range("A2").Value = "Rob Oudman"

# Update the active selection to A2
$ws.Range("A2").Select()
